# Update loading_percent results table (B2:O25) with the recalculated
# values for the "case with 380 kV" run. Columns F, J, L stay 0 (unchanged).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 24,14
$arr[0,0] = 7.969369857045184
$arr[0,1] = 5.195878383195836
$arr[0,2] = 5.977702694545851
$arr[0,3] = 12.81281460977439
$arr[0,4] = 0
$arr[0,5] = 26.58080220480687
$arr[0,6] = 13.73871556977121
$arr[0,7] = 19.70030335934489
$arr[0,8] = 0
$arr[0,9] = 7.96324129220353
$arr[0,10] = 0
$arr[0,11] = 13.30283403756381
$arr[0,12] = 18.08836794374145
$arr[0,13] = 20.62231329645772
$arr[1,0] = 7.666762284482881
$arr[1,1] = 5.019993534820838
$arr[1,2] = 5.857245263112068
$arr[1,3] = 12.60187461194797
$arr[1,4] = 0
$arr[1,5] = 26.5844931973599
$arr[1,6] = 13.77722697133975
$arr[1,7] = 19.77190954899463
$arr[1,8] = 0
$arr[1,9] = 7.677766151981976
$arr[1,10] = 0
$arr[1,11] = 13.12576387236318
$arr[1,12] = 18.14489452750833
$arr[1,13] = 20.67472979903178
$arr[2,0] = 7.475829547520485
$arr[2,1] = 4.907790184318357
$arr[2,2] = 5.783740505790169
$arr[2,3] = 12.47503073122984
$arr[2,4] = 0
$arr[2,5] = 26.59529441312765
$arr[2,6] = 13.80299161355505
$arr[2,7] = 19.81950432080768
$arr[2,8] = 0
$arr[2,9] = 7.495279580964001
$arr[2,10] = 0
$arr[2,11] = 13.01902840978592
$arr[2,12] = 18.18121419551629
$arr[2,13] = 20.71126297264361
$arr[3,0] = 7.396853865502632
$arr[3,1] = 4.86104869240199
$arr[3,2] = 5.753948618913691
$arr[3,3] = 12.42408768199937
$arr[3,4] = 0
$arr[3,5] = 26.60183916582953
$arr[3,6] = 13.81402350322689
$arr[3,7] = 19.83981107605815
$arr[3,8] = 0
$arr[3,9] = 7.426744033419517
$arr[3,10] = 0
$arr[3,11] = 12.97608529095523
$arr[3,12] = 18.19642133564918
$arr[3,13] = 20.7272421140579
$arr[4,0] = 7.383673196361558
$arr[4,1] = 4.853227114059469
$arr[4,2] = 5.749012864208463
$arr[4,3] = 12.4156758847452
$arr[4,4] = 0
$arr[4,5] = 26.60305523849489
$arr[4,6] = 13.81588750341833
$arr[4,7] = 19.84323802443049
$arr[4,8] = 0
$arr[4,9] = 7.418047756756295
$arr[4,10] = 0
$arr[4,11] = 12.96898951360837
$arr[4,12] = 18.19897106147404
$arr[4,13] = 20.72996130319375
$arr[5,0] = 7.474769017862632
$arr[5,1] = 4.907163875794791
$arr[5,2] = 5.783338003812564
$arr[5,3] = 12.47434057407622
$arr[5,4] = 0
$arr[5,5] = 26.59537400558042
$arr[5,6] = 13.80313823721993
$arr[5,7] = 19.81977449533238
$arr[5,8] = 0
$arr[5,9] = 7.494260093670199
$arr[5,10] = 0
$arr[5,11] = 13.01844695778231
$arr[5,12] = 18.18141763650895
$arr[5,13] = 20.71147405604959
$arr[6,0] = 7.866164280555027
$arr[6,1] = 5.136130030358618
$arr[6,2] = 5.936102243233372
$arr[6,3] = 12.73957249198682
$arr[6,4] = 0
$arr[6,5] = 26.58030212520118
$arr[6,6] = 13.75155453990029
$arr[6,7] = 19.72423960860751
$arr[6,8] = 0
$arr[6,9] = 7.866340882491975
$arr[6,10] = 0
$arr[6,11] = 13.24140003413413
$arr[6,12] = 18.10752444105029
$arr[6,13] = 20.6394826235004
$arr[7,0] = 8.588300386094859
$arr[7,1] = 5.550161607171153
$arr[7,2] = 6.237268791624063
$arr[7,3] = 13.27759524737951
$arr[7,4] = 0
$arr[7,5] = 26.61852634484894
$arr[7,6] = 13.66721355377164
$arr[7,7] = 19.56572305277156
$arr[7,8] = 0
$arr[7,9] = 8.536378473829831
$arr[7,10] = 0
$arr[7,11] = 13.6920071610683
$arr[7,12] = 17.97535575645713
$arr[7,13] = 20.53290395805399
$arr[8,0] = 9.08575603797301
$arr[8,1] = 5.831213250179655
$arr[8,2] = 6.456885280612151
$arr[8,3] = 13.6793544427299
$arr[8,4] = 0
$arr[8,5] = 26.68790105042817
$arr[8,6] = 13.61550564812409
$arr[8,7] = 19.46688668167289
$arr[8,8] = 0
$arr[8,9] = 8.989508414203819
$arr[8,10] = 0
$arr[8,11] = 14.02804989602822
$arr[8,12] = 17.88593592517299
$arr[8,13] = 20.47580217079859
$arr[9,0] = 9.303948697857544
$arr[9,1] = 5.953736116110319
$arr[9,2] = 6.555914963984646
$arr[9,3] = 13.86260758246794
$arr[9,4] = 0
$arr[9,5] = 26.72838355432262
$arr[9,6] = 13.59421097626828
$arr[9,7] = 19.42576189850044
$arr[9,8] = 0
$arr[9,9] = 9.18667761575492
$arr[9,10] = 0
$arr[9,11] = 14.18133550732549
$arr[9,12] = 17.84690817291007
$arr[9,13] = 20.45444886996778
$arr[10,0] = 9.385343444518698
$arr[10,1] = 5.999344353729148
$arr[10,2] = 6.593248568753738
$arr[10,3] = 13.9319969245829
$arr[10,4] = 0
$arr[10,5] = 26.7449896818351
$arr[10,6] = 13.58646766029726
$arr[10,7] = 19.41074165647643
$arr[10,8] = 0
$arr[10,9] = 9.260019332202869
$arr[10,10] = 0
$arr[10,11] = 14.23938601588976
$arr[10,12] = 17.83236536355465
$arr[10,13] = 20.44702895619836
$arr[11,0] = 9.367869279305596
$arr[11,1] = 5.989557195959465
$arr[11,2] = 6.585216153986948
$arr[11,3] = 13.91705406463887
$arr[11,4] = 0
$arr[11,5] = 26.74135662210272
$arr[11,6] = 13.58812106397622
$arr[11,7] = 19.41395192920071
$arr[11,8] = 0
$arr[11,9] = 9.244283157095424
$arr[11,10] = 0
$arr[11,11] = 14.22688446673341
$arr[11,12] = 17.83548693189751
$arr[11,13] = 20.44859732170545
$arr[12,0] = 9.310670086284121
$arr[12,1] = 5.957504253155848
$arr[12,2] = 6.55898999102568
$arr[12,3] = 13.86831680009119
$arr[12,4] = 0
$arr[12,5] = 26.72972421807464
$arr[12,6] = 13.59356750433604
$arr[12,7] = 19.4245150882158
$arr[12,8] = 0
$arr[12,9] = 9.192738151367328
$arr[12,10] = 0
$arr[12,11] = 14.18611157928769
$arr[12,12] = 17.84570700053125
$arr[12,13] = 20.45382507198043
$arr[13,0] = 9.275471958241125
$arr[13,1] = 5.937767603006744
$arr[13,2] = 6.542902812780616
$arr[13,3] = 13.83846102507421
$arr[13,4] = 0
$arr[13,5] = 26.7227650072126
$arr[13,6] = 13.59694535237387
$arr[13,7] = 19.43105735581905
$arr[13,8] = 0
$arr[13,9] = 9.160992278585326
$arr[13,10] = 0
$arr[13,11] = 14.1611359784432
$arr[13,12] = 17.85199780985497
$arr[13,13] = 20.45711400589675
$arr[14,0] = 9.07132783740415
$arr[14,1] = 5.823096894194084
$arr[14,2] = 6.450392295789246
$arr[14,3] = 13.66738180762546
$arr[14,4] = 0
$arr[14,5] = 26.6854344400804
$arr[14,6] = 13.61694214578938
$arr[14,7] = 19.46965157957294
$arr[14,8] = 0
$arr[14,9] = 8.976439654466285
$arr[14,10] = 0
$arr[14,11] = 14.01803595599832
$arr[14,12] = 17.88851958553601
$arr[14,13] = 20.4772908017017
$arr[15,0] = 8.943967740068235
$arr[15,1] = 5.751368504524579
$arr[15,2] = 6.393386102811983
$arr[15,3] = 13.5625016552302
$arr[15,4] = 0
$arr[15,5] = 26.66481456122732
$arr[15,6] = 13.62978019838843
$arr[15,7] = 19.49431129816292
$arr[15,8] = 0
$arr[15,9] = 8.860902045790162
$arr[15,10] = 0
$arr[15,11] = 13.93031588827537
$arr[15,12] = 17.9113462860968
$arr[15,13] = 20.49085362891715
$arr[16,0] = 8.869954229966831
$arr[16,1] = 5.709611777464983
$arr[16,2] = 6.36051759943755
$arr[16,3] = 13.50222792412635
$arr[16,4] = 0
$arr[16,5] = 26.65379501673065
$arr[16,6] = 13.63737395984508
$arr[16,7] = 19.50885594551557
$arr[16,8] = 0
$arr[16,9] = 8.793606003345298
$arr[16,10] = 0
$arr[16,11] = 13.87990422941944
$arr[16,12] = 17.92463094943077
$arr[16,13] = 20.49908963822776
$arr[17,0] = 8.844766233625984
$arr[17,1] = 5.69538843925595
$arr[17,2] = 6.349376427101958
$arr[17,3] = 13.481831291362
$arr[17,4] = 0
$arr[17,5] = 26.65020851933499
$arr[17,6] = 13.63998107876177
$arr[17,7] = 19.51384248505869
$arr[17,8] = 0
$arr[17,9] = 8.77067722050481
$arr[17,10] = 0
$arr[17,11] = 13.86284473859575
$arr[17,12] = 17.92915562266314
$arr[17,13] = 20.50195288025352
$arr[18,0] = 8.957604528381522
$arr[18,1] = 5.759056079479103
$arr[18,2] = 6.399463098379712
$arr[18,3] = 13.57366164696782
$arr[18,4] = 0
$arr[18,5] = 26.66692264183252
$arr[18,6] = 13.62839186594597
$arr[18,7] = 19.49164885542893
$arr[18,8] = 0
$arr[18,9] = 8.873288607412437
$arr[18,10] = 0
$arr[18,11] = 13.93964983639133
$arr[18,12] = 17.90890027468884
$arr[18,13] = 20.48936480805011
$arr[19,0] = 9.327504734158717
$arr[19,1] = 5.966940542078492
$arr[19,2] = 6.566698092133468
$arr[19,3] = 13.88263282871599
$arr[19,4] = 0
$arr[19,5] = 26.73310636354169
$arr[19,6] = 13.59195905260496
$arr[19,7] = 19.42139741980342
$arr[19,8] = 0
$arr[19,9] = 9.207914289225091
$arr[19,10] = 0
$arr[19,11] = 14.19808787651347
$arr[19,12] = 17.8426987180769
$arr[19,13] = 20.45227146670828
$arr[20,0] = 9.562061918565449
$arr[20,1] = 6.098200691761096
$arr[20,2] = 6.67500635503225
$arr[20,3] = 14.08450332378545
$arr[20,4] = 0
$arr[20,5] = 26.78379588388054
$arr[20,6] = 13.57001638525312
$arr[20,7] = 19.37870683740688
$arr[20,8] = 0
$arr[20,9] = 9.418893760617149
$arr[20,10] = 0
$arr[20,11] = 14.3669939855105
$arr[20,12] = 17.80080816144717
$arr[20,13] = 20.43191196136446
$arr[21,0] = 9.43755182696259
$arr[21,1] = 6.02857257964556
$arr[21,2] = 6.617303651655049
$arr[21,3] = 13.97679089200762
$arr[21,4] = 0
$arr[21,5] = 26.75606442832389
$arr[21,6] = 13.58155659007164
$arr[21,7] = 19.40119633999198
$arr[21,8] = 0
$arr[21,9] = 9.307006028685768
$arr[21,10] = 0
$arr[21,11] = 14.27686342998238
$arr[21,12] = 17.82304039203288
$arr[21,13] = 20.44242250049615
$arr[22,0] = 8.951441799800437
$arr[22,1] = 5.755582142100415
$arr[22,2] = 6.396715981094039
$arr[22,3] = 13.56861613647815
$arr[22,4] = 0
$arr[22,5] = 26.66596697680438
$arr[22,6] = 13.62901886758615
$arr[22,7] = 19.49285140129219
$arr[22,8] = 0
$arr[22,9] = 8.867691356237597
$arr[22,10] = 0
$arr[22,11] = 13.93542989258925
$arr[22,12] = 17.91000561394483
$arr[22,13] = 20.49003653797594
$arr[23,0] = 8.398407216521573
$arr[23,1] = 5.442101378114153
$arr[23,2] = 6.155901538747925
$arr[23,3] = 13.13059019726811
$arr[23,4] = 0
$arr[23,5] = 26.60092521461805
$arr[23,6] = 13.68822916129632
$arr[23,7] = 19.60551446656448
$arr[23,8] = 0
$arr[23,9] = 8.361811572869764
$arr[23,10] = 0
$arr[23,11] = 13.56900112857131
$arr[23,12] = 18.00975570305352
$arr[23,13] = 20.53290395805399
$ws.Range("B2:O25").Value = $arr
Write-Output "done"
